$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$errZhCn = "Handback file name: tutq1y5d.nxw is different with handoff file name: f3903bd5-f096-4661-a889-b1aab169d0ba.c8d1622f5a38209faee6e8239b4ae03280ca43b2.zh-cn."
$errDeDe = "Handback file name: tutq1y5d.nxw is different with handoff file name: f3903bd5-f096-4661-a889-b1aab169d0ba.c8d1622f5a38209faee6e8239b4ae03280ca43b2.de-de."

# Overview sheet: row for f3903bd5 file (row 3) - zh-cn (E3) and de-de (F3) status columns
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn sheet: status column (C3) + error detail column (P3) for f3903bd5 row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = $errZhCn
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1

# de-de sheet: status column (C3) + error detail column (P3) for f3903bd5 row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = $errDeDe
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1
